# Variable Definitions RQ8.xlsx update:
# "updated RQ5 and variables to base off text instead of delegation"
#
# 1) Rename Delegator_ID -> Participant_ID (and related text) throughout.
# 2) Re-wire the "Source" row for TA1_Name/Attribute to the response log
#    (was pointing at the loading log).
# 3) Remove a duplicated "evac" column that had crept into each of the
#    8 Patient{n} column groups (each group was 7 cols: time, order, evac,
#    assess, treat, evac(dup), tag -- the dup is removed, leaving 6 cols).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------

$ws.Range("B1").Value = "Participant_ID"
$ws.Range("G1").Value = "Alignment score (Participant|selected target)"

$ws.Range("C2").Value = "Text scenario response log"
$ws.Range("D2").Value = "Text scenario response log"

# Row 2 "Source" column text (B/E) and row 3 "Definition" text are also
# rewritten even where the value itself doesn't move columns.
$ws.Range("B2").Value = "Text scenario response log"
$ws.Range("E2").Value = "Text scenario loading log"

$ws.Range("B3").Value = "Used to track and identify participants, also called delegator ID"
$ws.Range("E3").Value = "Scenario presented to the participant in text scenarios"
$ws.Range("G3").Value = "Calculated alignment score between the KDMA measurement of a participant and a selected target"

# --- Structural fix: drop the duplicated "evac" column in each of the  ---
# --- 8 Patient{n} blocks (columns V, AC, AJ, AQ, AX, BE, BL, BS).      ---
# Deleted right-to-left so earlier deletions don't shift later targets.

$ws.Range("BS1:BS4").EntireColumn.Delete()
$ws.Range("BL1:BL4").EntireColumn.Delete()
$ws.Range("BE1:BE4").EntireColumn.Delete()
$ws.Range("AX1:AX4").EntireColumn.Delete()
$ws.Range("AQ1:AQ4").EntireColumn.Delete()
$ws.Range("AJ1:AJ4").EntireColumn.Delete()
$ws.Range("AC1:AC4").EntireColumn.Delete()
$ws.Range("V1:V4").EntireColumn.Delete()

# --- View state -------------------------------------------------------
[void]$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 4
[void]$ws.Range("H3").Select()
